# Update cryptos price (D) and 1h volume change (E) columns with refreshed
# values. Numeric-looking Price values are written with a leading apostrophe
# so Excel stores them as text (matching the sheet's existing inline-string
# cells) instead of auto-converting them to numbers; the style is then reset
# to Normal so no stray "quote prefix" formatting is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.148.90"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.247.48"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'247.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("D7").Value = "'74.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.35%  "
$ws.Range("D10").Value = "'41.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("D11").Value = "'0.0941"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.13%  "
$ws.Range("D12").Value = "'7.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.19%  "
$ws.Range("E13").Value = "  -4.49%  "
$ws.Range("D14").Value = "2.583.30"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "'14.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "2.251.66"
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").Value = "42.069.16"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "0.0₃0979"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "'71.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "'2.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.46%  "
$ws.Range("D23").Value = "'230.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D25").Value = "'11.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("E26").Value = "  -7.97%  "
$ws.Range("D27").Value = "'7.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +24.71%  "
$ws.Range("D28").Value = "'2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("D29").Value = "'171.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").Value = "'2.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.91%  "
$ws.Range("D31").Value = "'20.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "'0.0825"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("E33").Value = "  -5.72%  "
$ws.Range("D34").Value = "'30.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("D35").Value = "'0.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("D36").Value = "'4.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").Value = "'4.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").Value = "'13.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "'2.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").Value = "'62.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").Value = "'108.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("D44").Value = "'0.202"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("E45").Value = "  -4.70%  "
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  -3.72%  "
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "'2.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("E51").Value = "  -0.92%  "
